# Deploy updated output folder
# Update the "Metadata" worksheet's Title and Date values to reflect the
# newly generated IG output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 5: Title -> new value
$ws.Range("B5").Value = "NG-Imm Vaccine Site VS"

# Row 8: Date -> new generation timestamp
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
